{"js": "const paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nconst oldText = \"\u041d\u0430 \u0441\u0445\u0435\u043c\u0443 \u0431\u044b\u043b \u0434\u043e\u0431\u0430\u0432\u043b\u0435\u043d \u0438\u0441\u0442\u043e\u0447\u043d\u0438\u043a \u043d\u0430\u043f\u0440\u044f\u0436\u0435\u043d\u0438\u044f. \";\nlet newText = \"\u041d\u0430 \u0441\u0445\u0435\u043c\u0443 \u0431\u044b\u043b \u0434\u043e\u0431\u0430\u0432\u043b\u0435\u043d \u0430\u043c\u043f\u0435\u0440\u043c\u0435\u0442\u0440. \";\nfor (let i = 0; i < 10; i++) {\n  newText += \"\u041d\u0430 \u0441\u0445\u0435\u043c\u0443 \u0431\u044b\u043b \u0434\u043e\u0431\u0430\u0432\u043b\u0435\u043d \u0432\u043e\u043b\u044c\u0442\u043c\u0435\u0442\u0440. \";\n}\n\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  const p = paragraphs.items[i];\n  p.load(\"text\");\n}\nawait context.sync();\n\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  const p = paragraphs.items[i];\n  if (p.text === oldText) {\n    p.insertText(newText, Word.InsertLocation.replace);\n  }\n}\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n$oldText = \"\u041d\u0430 \u0441\u0445\u0435\u043c\u0443 \u0431\u044b\u043b \u0434\u043e\u0431\u0430\u0432\u043b\u0435\u043d \u0438\u0441\u0442\u043e\u0447\u043d\u0438\u043a \u043d\u0430\u043f\u0440\u044f\u0436\u0435\u043d\u0438\u044f. \"\n$newText = \"\u041d\u0430 \u0441\u0445\u0435\u043c\u0443 \u0431\u044b\u043b \u0434\u043e\u0431\u0430\u0432\u043b\u0435\u043d \u0430\u043c\u043f\u0435\u0440\u043c\u0435\u0442\u0440. \"\nfor ($i = 0; $i -lt 10; $i++) {\n    $newText += \"\u041d\u0430 \u0441\u0445\u0435\u043c\u0443 \u0431\u044b\u043b \u0434\u043e\u0431\u0430\u0432\u043b\u0435\u043d \u0432\u043e\u043b\u044c\u0442\u043c\u0435\u0442\u0440. \"\n}\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$found = $find.Execute(\n    $oldText,\n    $false,\n    $false,\n    $false,\n    $false,\n    $false,\n    $true,\n    1,\n    $false,\n    $newText,\n    2\n)\n"}
